$d = $word.ActiveDocument

$rng = $d.Content.Duplicate
$found = $rng.Find.Execute("https://www.youtube.com/watch?v=Ic10mDeUt-8")
Write-Host "Find result: $found"

try {
  $rng.Font.Color = "wdThemeColorAccent3"
  Write-Host "Set string OK"
} catch {
  Write-Host "Set string FAILED: $_"
}
Write-Host "Color now: $($rng.Font.Color)"

try {
  $x = $rng.Font.TextColor
  Write-Host "Font.TextColor = $x"
} catch {
  Write-Host "Get Font.TextColor FAILED: $_"
}

try {
  $rng.Font.TextColor.ObjectThemeColor = 4
  Write-Host "Set Font.TextColor.ObjectThemeColor OK"
} catch {
  Write-Host "Set Font.TextColor.ObjectThemeColor FAILED: $_"
}

try {
  $rng.Font.ObjectThemeColor = 4
  Write-Host "Set Font.ObjectThemeColor OK"
} catch {
  Write-Host "Set Font.ObjectThemeColor FAILED: $_"
}
